$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("High Elf Wrong Profiles")
$ws.Activate()

# Row 2 - Eltharion (main model)
$ws.Range("F2").Value = 467
$ws.Range("G2").Value = 301
$ws.Range("H2").Formula = "=F2-G2"

# Row 4 - Tyrion (main model)
$ws.Range("F4").Value = 425
$ws.Range("G4").Value = 250
$ws.Range("H4").Formula = "=F4-G4"

# Row 6 - Drachenprinzen von Caledor
$ws.Range("F6").Value = 43
$ws.Range("G6").Value = 10

# Row 12 - Ellyrianische Grenzreiter
$ws.Range("F12").Value = 25
$ws.Range("G12").Value = 2

# Row 18 - Schwertmeister von Hoeth
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 14

# Row 21 - Silberhelme
$ws.Range("F21").Value = 31
$ws.Range("G21").Value = 4

# Row 27 - Weiße Löwen von Chrace
$ws.Range("F27").Value = 16
$ws.Range("G27").Value = 14

# Row 30 - Zwillingskämpfer
$ws.Range("F30").Value = 9
$ws.Range("G30").Value = 11

# Row 33 - Tiranoc Streitwagen mit einem Elfen
$ws.Range("F33").Value = 72
$ws.Range("G33").Value = 2

# Row 36 - Tiranoc Streitwagen mit zwei Elfen
$ws.Range("F36").Value = 84
$ws.Range("G36").Value = 2

# Row 39 - Repetier-Speerschleuder
$ws.Range("F39").Value = 100
$ws.Range("G39").Value = 2

# Restore the selection to match where the author left off
$ws.Range("F6").Select()
